$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Footer "date last figured out" placeholder: 1/29/2022 -> 2/2/2022
#    This text lives on the slide master and every custom (slide) layout,
#    inside the shape named "Date Placeholder N". Update it everywhere it
#    appears so the whole deck's footer date is refreshed for Week 5.
# ---------------------------------------------------------------------------
$newDate = "2/2/2022"

function Update-DatePlaceholder($container) {
    for ($i = 1; $i -le $container.Shapes.Count; $i++) {
        $shp = $container.Shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

Update-DatePlaceholder $p.SlideMaster

for ($li = 1; $li -le $p.SlideMaster.CustomLayouts.Count; $li++) {
    Update-DatePlaceholder $p.SlideMaster.CustomLayouts.Item($li)
}

# ---------------------------------------------------------------------------
# 2) Reposition the small GitHub icon picture on slide 8 ("GitHub: Pull").
#    Old offset: x=3455582 y=3429000 (EMU)  ->  New offset: x=3395947 y=4001294 (EMU)
#    PowerPoint's Shape.Left/Top are expressed in points (1 pt = 12700 EMU);
#    the literal constants below are chosen so the round trip through the
#    host's single-precision storage reproduces the exact target EMU values.
# ---------------------------------------------------------------------------
$slide8 = $p.Slides.Item(8)
for ($i = 1; $i -le $slide8.Shapes.Count; $i++) {
    $shp = $slide8.Shapes.Item($i)
    if ($shp.Name -eq "Picture 5") {
        $shp.Left = 267.3974304199
        $shp.Top = 315.0625305176
    }
}
